# [ADDITIONAL SCRAPING] Add a "Player Info" sheet ahead of "ODI Batting" and
# extend the ODI Batting data with the extra scraped match-code / 2nd innings
# information.

$wb = $excel.ActiveWorkbook

# --- New "Player Info" sheet (inserted as the first/left-most sheet) -------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row.
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header styling used on the other
# sheet's header row.
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Data row - ID is numeric-looking text, so force text formatting before
# assigning the value so it is stored as a string, not a number.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6150"
$playerInfo.Range("B2").Value = "Rachin Ravindra"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- "ODI Batting" sheet updates -------------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")

# MATCH_CARD_LINK -> MATCH_CODE, and store just the bare code instead of the
# full scorecard URL.
$odi.Range("D1").Value = "MATCH_CODE"
$odi.Range("D2").NumberFormat = "@"
$odi.Range("D2").Value = "4735"

# Newly scraped second-innings row (player did not bat).
$odi.Range("A3").NumberFormat = "@"
$odi.Range("A3").Value = "2"
$odi.Range("C3").Value = "31/03/2023"
$odi.Range("D3").NumberFormat = "@"
$odi.Range("D3").Value = "4745"
$odi.Range("E3").Value = "2nd"
$odi.Range("F3").Value = "Sri Lanka"
$odi.Range("G3").Value = "Seddon Park"
$odi.Range("H3").Value = "did not bat"
$odi.Range("I3").Value = "-"
$odi.Range("J3").Value = "-"

# Make "Player Info" the active sheet/selection, matching a freshly-added
# sheet being the one left in focus.
$playerInfo.Activate()
$playerInfo.Range("A1").Select() | Out-Null
